$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" column header
$ws.Range("T1").Value = "Total"

# Row totals for the existing disease-group rows
$ws.Range("T2").Value = 80152
$ws.Range("T3").Value = 7774
$ws.Range("T4").Value = 34410
$ws.Range("T5").Value = 13065
$ws.Range("T6").Value = 48871

# New row 7: "Outros"
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 7242
$ws.Range("C7").Value = 420
$ws.Range("D7").Value = 560
$ws.Range("E7").Value = 2440
$ws.Range("F7").Value = 3513
$ws.Range("G7").Value = 3568
$ws.Range("H7").Value = 3884
$ws.Range("I7").Value = 4037
$ws.Range("J7").Value = 4556
$ws.Range("K7").Value = 4958
$ws.Range("L7").Value = 5165
$ws.Range("M7").Value = 5431
$ws.Range("N7").Value = 5061
$ws.Range("O7").Value = 5101
$ws.Range("P7").Value = 5306
$ws.Range("Q7").Value = 6121
$ws.Range("R7").Value = 18085
$ws.Range("S7").Value = 712
$ws.Range("T7").Value = 86160

# New row 8: "Total"
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 8150
$ws.Range("C8").Value = 651
$ws.Range("D8").Value = 822
$ws.Range("E8").Value = 2988
$ws.Range("F8").Value = 4248
$ws.Range("G8").Value = 4657
$ws.Range("H8").Value = 5583
$ws.Range("I8").Value = 6605
$ws.Range("J8").Value = 8796
$ws.Range("K8").Value = 11955
$ws.Range("L8").Value = 15755
$ws.Range("M8").Value = 19908
$ws.Range("N8").Value = 22099
$ws.Range("O8").Value = 24159
$ws.Range("P8").Value = 26547
$ws.Range("Q8").Value = 29977
$ws.Range("R8").Value = 76627
$ws.Range("S8").Value = 905
$ws.Range("T8").Value = 270432
